$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Unhide columns D:O (they were hidden, now shown) ---
$ws.Range("D1:O1").EntireColumn.Hidden = $false

# --- Update data cells ---
# NOTE: rows 28/35/42/49/56/77/85/90 are hidden rows. Writing a cell value
# into a hidden row causes this host's best-effort row-autofit to stamp a
# spurious ht/customHeight on the row. Briefly unhiding the row around the
# writes avoids that, then we restore the original hidden state.

function Set-RowValues($RowNum, $Values, $NumberFormats) {
    $row = $ws.Rows($RowNum)
    $wasHidden = $row.Hidden
    $row.Hidden = $false

    foreach ($addr in $Values.Keys) {
        $ws.Range($addr).Value = $Values[$addr]
    }
    foreach ($addr in $NumberFormats.Keys) {
        $ws.Range($addr).NumberFormat = $NumberFormats[$addr]
    }

    $row.Hidden = $wasHidden
}

# Row 28: D/E already use the "#,##0" style (s=15) - just change values.
Set-RowValues 28 @{
    "D28" = 732491
    "E28" = 1035114
} @{}

# Row 35: D/E move from General (s=12) to "#,##0" (s=15).
Set-RowValues 35 @{
    "D35" = 824391
    "E35" = 742712
} @{
    "D35" = "#,##0"
    "E35" = "#,##0"
}

# Row 42
Set-RowValues 42 @{
    "D42" = 779804
    "E42" = 786990
} @{
    "D42" = "#,##0"
    "E42" = "#,##0"
}

# Row 49
Set-RowValues 49 @{
    "D49" = 836380
    "E49" = 795511
} @{
    "D49" = "#,##0"
    "E49" = "#,##0"
}

# Row 56
Set-RowValues 56 @{
    "D56" = 663218
    "E56" = 771495
} @{
    "D56" = "#,##0"
    "E56" = "#,##0"
}

# Row 77
Set-RowValues 77 @{
    "D77" = 714382
    "E77" = 969381
} @{
    "D77" = "#,##0"
    "E77" = "#,##0"
}

# Row 85
Set-RowValues 85 @{
    "D85" = 766641
    "E85" = 648864
} @{
    "D85" = "#,##0"
    "E85" = "#,##0"
}

# Row 90: fill in previously-empty KPI cells (P, Q, R, T, U).
Set-RowValues 90 @{
    "P90" = 4
    "Q90" = 4
    "R90" = 29244.63
    "T90" = 0
    "U90" = 1
} @{
    "U90" = "0%"
}
